# B6-PowerPoint.pptx edit: Tue, Apr 28, 2020  6:06:07 AM
#
# 1) Three tables (slides 14, 15, 16 - each table is the first shape on its
#    slide) get their table style switched from the default
#    {1C583B64-3332-420B-8458-AA28EB2CCBFE} to
#    {1076CE3C-3719-4F98-8C36-683E877EB9FB}.
# 2) The deck's theme is recoloured from the "Integral" / "Red Violet"
#    palette to the plain "Office" palette (the 12 DrawingML theme colors
#    that back the slide master's theme part).

$p = $ppt.ActivePresentation

# --- 1) Table styles -------------------------------------------------------

$newStyleId = "{1076CE3C-3719-4F98-8C36-683E877EB9FB}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    $shape = $slide.Shapes.Item(1)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle($newStyleId)
    }
}

# --- 2) Theme colors --------------------------------------------------------

function Get-RGBFromHex($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink - the "Office" theme palette.
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = Get-RGBFromHex($officeThemeColors[$i - 1])
}
